$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the "F1" / "F2" cross labels into sub-groups (F1a/F1b, F2a/F2b) and
# move "P2" up next to "P1" in the shared-strings table by rewriting the
# cross column (A3:A7) with the new labels.
$ws.Range("A3").Value = "F1a"
$ws.Range("A4").Value = "F1b"
$ws.Range("A5").Value = "F2a"
$ws.Range("A6").Value = "F2b"
$ws.Range("A7").Value = "P2"

# Update the sheet's current selection to D18.
[void]$ws.Range("D18").Select()
